$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Z (26) - shifts COLL_SYMBOL..RESTRAINTS_DISCOUNT
# one column to the right and makes room for the new "STAT" column.
$ws.Columns("Z:Z").Insert()

# Header for the newly inserted column.
$ws.Range("Z1").Value = "STAT"

# Data value for row 2 in the new column.
$ws.Range("Z2").Value = "K"

# Give the inserted column its own (narrower) width.
$ws.Columns("Z:Z").ColumnWidth = 5.7

# Update the view: scrolled so column X is the first visible column, and the
# active selection moved to AE9.
$ws.Application.ActiveWindow.ScrollColumn = 24
$ws.Range("AE9").Select()
